$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cD = $ws.Range("D2")
$cD.NumberFormat = "@"
$cD.Value = '66.992.31'
$cD.Style = "Normal"
$ws.Range("E2").Value = '  +2.03%  '

$cD = $ws.Range("D3")
$cD.NumberFormat = "@"
$cD.Value = '3.281.13'
$cD.Style = "Normal"
$ws.Range("E3").Value = '  -0.81%  '

$cD = $ws.Range("D4")
$cD.NumberFormat = "@"
$cD.Value = '0.998'
$cD.Style = "Normal"
$ws.Range("E4").Value = '  -0.12%  '

$cD = $ws.Range("D5")
$cD.NumberFormat = "@"
$cD.Value = '569.98'
$cD.Style = "Normal"
$ws.Range("E5").Value = '  -1.61%  '

$cD = $ws.Range("D6")
$cD.NumberFormat = "@"
$cD.Value = '174.98'
$cD.Style = "Normal"
$ws.Range("E6").Value = '  -4.81%  '

$cD = $ws.Range("D7")
$cD.NumberFormat = "@"
$cD.Value = '0.999'
$cD.Style = "Normal"
$ws.Range("E7").Value = '  -0.10%  '

$cD = $ws.Range("D9")
$cD.NumberFormat = "@"
$cD.Value = '3.279.57'
$cD.Style = "Normal"
$ws.Range("E9").Value = '  -0.71%  '

$cD = $ws.Range("D10")
$cD.NumberFormat = "@"
$cD.Value = '0.174'
$cD.Style = "Normal"
$ws.Range("E10").Value = '  -2.46%  '

$ws.Range("E11").Value = '  -0.14%  '

$cD = $ws.Range("D12")
$cD.NumberFormat = "@"
$cD.Value = '45.50'
$cD.Style = "Normal"
$ws.Range("E12").Value = '  -2.51%  '

$ws.Range("E13").Value = '  +1.17%  '

$cD = $ws.Range("D14")
$cD.NumberFormat = "@"
$cD.Value = '688.34'
$cD.Style = "Normal"

$cD = $ws.Range("D15")
$cD.NumberFormat = "@"
$cD.Value = '3.805.50'
$cD.Style = "Normal"
$ws.Range("E15").Value = '  -0.80%  '

$cD = $ws.Range("D16")
$cD.NumberFormat = "@"
$cD.Value = '8.28'
$cD.Style = "Normal"
$ws.Range("E16").Value = '  -2.11%  '

$cD = $ws.Range("D17")
$cD.NumberFormat = "@"
$cD.Value = '67.022.03'
$cD.Style = "Normal"

$ws.Range("E18").Value = '  +0.96%  '

$cD = $ws.Range("D19")
$cD.NumberFormat = "@"
$cD.Value = '3.283.22'
$cD.Style = "Normal"
$ws.Range("E19").Value = '  -0.70%  '

$cD = $ws.Range("D20")
$cD.NumberFormat = "@"
$cD.Value = '17.30'
$cD.Style = "Normal"
$ws.Range("E20").Value = '  -3.02%  '

$cD = $ws.Range("D21")
$cD.NumberFormat = "@"
$cD.Value = '10.71'
$cD.Style = "Normal"
$ws.Range("E21").Value = '  -2.60%  '

$cD = $ws.Range("D22")
$cD.NumberFormat = "@"
$cD.Value = '0.885'
$cD.Style = "Normal"
$ws.Range("E22").Value = '  -0.56%  '

$ws.Range("E23").Value = '  -4.63%  '

$cD = $ws.Range("D24")
$cD.NumberFormat = "@"
$cD.Value = '5.11'
$cD.Style = "Normal"
$ws.Range("E24").Value = '  +1.84%  '

$cD = $ws.Range("D25")
$cD.NumberFormat = "@"
$cD.Value = '99.12'
$cD.Style = "Normal"
$ws.Range("E25").Value = '  -1.08%  '

$ws.Range("E26").Value = '  -2.50%  '

$ws.Range("E27").Value = '  -1.34%  '

$ws.Range("E28").Value = '  -1.48%  '

$cD = $ws.Range("D29")
$cD.NumberFormat = "@"
$cD.Value = '32.94'
$cD.Style = "Normal"
$ws.Range("E29").Value = '  +6.67%  '

$cD = $ws.Range("D30")
$cD.NumberFormat = "@"
$cD.Value = '8.37'
$cD.Style = "Normal"
$ws.Range("E30").Value = '  -0.02%  '

$cD = $ws.Range("D31")
$cD.NumberFormat = "@"
$cD.Value = '6.74'
$cD.Style = "Normal"
$ws.Range("E31").Value = '  +2.63%  '

$cD = $ws.Range("D32")
$cD.NumberFormat = "@"
$cD.Value = '572.51'
$cD.Style = "Normal"
$ws.Range("E32").Value = '  -3.66%  '

$cD = $ws.Range("D33")
$cD.NumberFormat = "@"
$cD.Value = '3.881.48'
$cD.Style = "Normal"
$ws.Range("E33").Value = '  +0.66%  '

$cD = $ws.Range("D34")
$cD.NumberFormat = "@"
$cD.Value = '10.79'
$cD.Style = "Normal"
$ws.Range("E34").Value = '  -1.34%  '

$ws.Range("E35").Value = '  -2.16%  '

$cD = $ws.Range("D36")
$cD.NumberFormat = "@"
$cD.Value = '1.00'
$cD.Style = "Normal"
$ws.Range("E36").Value = '  +0.00%  '

$cD = $ws.Range("D37")
$cD.NumberFormat = "@"
$cD.Value = '55.49'
$cD.Style = "Normal"
$ws.Range("E37").Value = '  -0.28%  '

$cD = $ws.Range("D38")
$cD.NumberFormat = "@"
$cD.Value = '3.30'
$cD.Style = "Normal"
$ws.Range("E38").Value = '  -11.94%  '

$ws.Range("E39").Value = '  +1.02%  '

$ws.Range("E40").Value = '  -0.36%  '

$cD = $ws.Range("D41")
$cD.NumberFormat = "@"
$cD.Value = '3.35'
$cD.Style = "Normal"
$ws.Range("E41").Value = '  -1.81%  '

$cD = $ws.Range("D42")
$cD.NumberFormat = "@"
$cD.Value = '31.70'
$cD.Style = "Normal"
$ws.Range("E42").Value = '  -2.00%  '

$cD = $ws.Range("D43")
$cD.NumberFormat = "@"
$cD.Value = '0.0₃0669'
$cD.Style = "Normal"
$ws.Range("E43").Value = '  -4.44%  '

$ws.Range("E44").Value = '  -4.44%  '

$cD = $ws.Range("D45")
$cD.NumberFormat = "@"
$cD.Value = '0.327'
$cD.Style = "Normal"
$ws.Range("E45").Value = '  -1.84%  '

$ws.Range("E46").Value = '  -0.86%  '

$ws.Range("E47").Value = '  +0.31%  '

$ws.Range("E49").Value = '  +7.46%  '

$ws.Range("E50").Value = '  +0.31%  '

$cD = $ws.Range("D51")
$cD.NumberFormat = "@"
$cD.Value = '130.46'
$cD.Style = "Normal"
$ws.Range("E51").Value = '  -0.21%  '
